$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 1004:1005 (shifts existing rows 1004+ down by 2)
$ws.Rows("1004:1005").Insert()

# Populate new row 1004
$ws.Range("A1004").Value = 9
$ws.Range("B1004").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1004").Value = "Metropolitana"
$ws.Range("D1004").Value = 45147
$ws.Range("E1004").Value = 13
$ws.Range("F1004").Value = 100112040
$ws.Range("G1004").Value = "Cilantro"
$ws.Range("H1004").Value = "Sin especificar"
$ws.Range("I1004").Value = "Primera"
$ws.Range("J1004").Value = 70
$ws.Range("K1004").Value = 6000
$ws.Range("L1004").Value = 6000
$ws.Range("M1004").Value = 6000
$ws.Range("N1004").Value = "`$/caja 36 atados"
$ws.Range("O1004").Value = "Región Metropolitana"
$ws.Range("P1004").Value = 167
$ws.Range("Q1004").Value = 36
$ws.Range("R1004").Value = "Hortaliza"

# Populate new row 1005
$ws.Range("A1005").Value = 9
$ws.Range("B1005").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1005").Value = "Metropolitana"
$ws.Range("D1005").Value = 45147
$ws.Range("E1005").Value = 13
$ws.Range("F1005").Value = 100112040
$ws.Range("G1005").Value = "Cilantro"
$ws.Range("H1005").Value = "Sin especificar"
$ws.Range("I1005").Value = "Primera"
$ws.Range("J1005").Value = 160
$ws.Range("K1005").Value = 9000
$ws.Range("L1005").Value = 11000
$ws.Range("M1005").Value = 10000
$ws.Range("N1005").Value = "`$/docena de atados"
$ws.Range("O1005").Value = "Región Metropolitana"
$ws.Range("P1005").Value = 3333
$ws.Range("Q1005").Value = 3
$ws.Range("R1005").Value = "Hortaliza"
